$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.779.09"
$ws.Range("D3").Value = "1.565.07"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'206.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "'0.489"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'21.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").Value = "'0.0860"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.784.28"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.552.83"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("D17").Value = "26.775.66"
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'7.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'214.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'9.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "'152.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "1.385.86"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "'0.921"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.38%  "
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").Value = "'0.521"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "'0.814"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'0.992"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "'63.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "1.698.65"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'85.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "0.0₇0985"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'0.0949"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("E51").Value = "  -0.69%  "
